$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvoiceLog")

# --- Grow Table1 by one column ("Column1"), 28 -> 29 columns --------------
$table = $ws.ListObjects.Item("Table1")
$table.ListColumns.Add() | Out-Null
$newColIndex = $table.ListColumns.Count
$ws.Range("AC1").Value = "Column1"

# Give the new column's data cells the same plain centred look used by the
# other un-formatted columns in the table (style index 1: center/center).
$dataRng = $ws.Range("AC2:AC14")
$dataRng.HorizontalAlignment = -4108
$dataRng.VerticalAlignment = -4108

# --- Populate the first saved-invoice row (row 2) with sample data --------
$ws.Range("B2").Value = 45000
$ws.Range("C2").Value = "TechZenith Team"
$ws.Range("D2").Value = "Corporetior"
$ws.Range("E2").Value = "122 Innovation prive, silicon Valley, Cit 95054"
$ws.Range("F2").Value = "438a' Street, Anytown, CA 12345"
$ws.Range("J2").Value = "TZS-2023-0001"
$ws.Range("K2").Value = 199549990
$ws.Range("L2").Value = 1234567
$ws.Range("P2").Value = "Cand ServerhostingUINcth) Website Designs Dav satware Custanization Searty Aadt"
$ws.Range("R2").Value = "2 x1 x 5 x 1"
$ws.Range("U2").Value = "R500.00 R 3000.00 R150.00 2500.00"
$ws.Range("V2").Value = "21000.00 R3000.00 R750.00 2500.00"
$ws.Range("AA2").Value = "R5250.00"
$ws.Range("AB2").Value = "R525250.00"
$ws.Range("AC2").Value = "12-3456789"
